# The sheet lists reporting-organisation groups with columns:
#   A = code, B = status, C = group-name, D = group-code
# The upstream codelist source swapped the column order so that the
# "code" value comes before the "name" value (group-code then group-name).
# Concretely, for every row (including the header row) the contents of
# column C and column D are exchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last populated row in column A (94 rows incl. header).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 1; $r -le $lastRow; $r++) {
    $cCell = $ws.Cells.Item($r, 3)
    $dCell = $ws.Cells.Item($r, 4)

    $cVal = $cCell.Value2
    $dVal = $dCell.Value2

    $cCell.Value2 = $dVal
    $dCell.Value2 = $cVal
}
